$d = $word.ActiveDocument

# --- Edit 1: Replace "Dump should have been imported..." paragraph with 3 new paragraphs ---
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Dump should have been imported from /files/caTissue/dump. ")
if (-not $found1) { Write-Output "ERROR: edit1 text not found" }
$rng1.Text = ""
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Import latest dump located at </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:spacing w:after="0"/></w:pPr><w:r><w:t>Oracle: https://ncisvn.nci.nih.gov/svn/catissue_persistent/caTissue Database Dump/v2.0/Oracle</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:spacing w:after="0"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>MySQL</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: https://ncisvn.nci.nih.gov/svn/catissue_persistent/caTissue Database Dump/v2.0/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MySQL</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and deploy application.</w:t></w:r></w:p>'
$rng1.InsertXML($xml1)

# --- Edit 2: Remove lastRenderedPageBreak before "3) " (Expected Output section) ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("3) Edit Participant")
if (-not $found2) { Write-Output "ERROR: edit2 text not found" }
$prng2 = $rng2.Paragraphs(1).Range
$prng2.Text = ""
$xml2 = '<w:p w:rsidR="00C121B9" w:rsidRPr="00D2500B" w:rsidRDefault="00C121B9" w:rsidP="00C121B9"><w:pPr><w:spacing w:after="0"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00D2500B"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">3) </w:t></w:r><w:r w:rsidRPr="00D2500B"><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Edit Participant</w:t></w:r><w:r w:rsidRPr="00D2500B"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> page is displayed on the RHS and LHS </w:t></w:r><w:r w:rsidRPr="00D2500B"><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Specimen Details</w:t></w:r><w:r w:rsidRPr="00D2500B"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> section should auto populate the 2 event points as </w:t></w:r></w:p>'
$prng2.InsertXML($xml2)

# --- Edit 3: Add lastRenderedPageBreak before "Aliquot 1" in table 2, row 5 ---
$t2 = $d.Tables.Item(2)
$cell3 = $t2.Cell(5,1)
$prng3 = $cell3.Range.Paragraphs(1).Range
$prng3.Text = ""
$xml3 = '<w:p w:rsidR="00C121B9" w:rsidRPr="00D2500B" w:rsidRDefault="00C121B9" w:rsidP="00853A5B"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00D2500B"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Aliquot 1</w:t></w:r></w:p>'
$prng3.InsertXML($xml3)

# --- Edit 4: Remove lastRenderedPageBreak before "Aliquot 4" in table 2, row 8 ---
$cell4 = $t2.Cell(8,1)
$prng4 = $cell4.Range.Paragraphs(1).Range
$prng4.Text = ""
$xml4 = '<w:p w:rsidR="00C121B9" w:rsidRPr="00D2500B" w:rsidRDefault="00C121B9" w:rsidP="00853A5B"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00D2500B"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Aliquot 4</w:t></w:r></w:p>'
$prng4.InsertXML($xml4)

# --- Edit 5: Add lastRenderedPageBreak before "parent_id" run ---
$rng5 = $d.Content
$found5 = $rng5.Find.Execute("Containment or reference type objects getting added will have a ")
if (-not $found5) { Write-Output "ERROR: edit5 text not found" }
$prng5 = $rng5.Paragraphs(1).Range
$prng5.Text = ""
$xml5 = '<w:p w:rsidR="005862A2" w:rsidRDefault="00C121B9" w:rsidP="005862A2"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="005862A2"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Object_ID</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="005862A2"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> is the unique ID of the object inserted. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="005862A2"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Parent_ID</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="005862A2"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> will be null for the main object (Specimen). Containment or reference type objects getting added will have a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="005862A2"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>parent_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="005862A2"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> equal to the ID of the main Object being inserted. This table refers to CATISSUE_AUDIT_EVENT_LOG table which relates to the CATISSUE_AUDIT_EVENT table.</w:t></w:r></w:p>'
$prng5.InsertXML($xml5)

# --- Edit 6: Remove lastRenderedPageBreak before "In CATISSUE_AUDIT_EVENT_DETAILS..." ---
$rng6 = $d.Content
$found6 = $rng6.Find.Execute("In CATISSUE_AUDIT_EVENT_DETAILS")
if (-not $found6) { Write-Output "ERROR: edit6 text not found" }
$prng6 = $rng6.Paragraphs(1).Range
$prng6.Text = ""
$xml6 = '<w:p w:rsidR="005862A2" w:rsidRDefault="00C121B9" w:rsidP="005862A2"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="005862A2"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>In CATISSUE_AUDIT_EVENT_DETAILS table Element name contains the list of attributes that are in CATISSUE_SPECIMEN.ID of all the reference and containment association classes should also be audited.</w:t></w:r></w:p>'
$prng6.InsertXML($xml6)

Write-Output "ALL DONE"
